$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.000.52"
$ws.Range("E2").Value = "  -1.71%  "

$ws.Range("D3").Value = "2.300.84"
$ws.Range("E3").Value = "  -1.96%  "

$ws.Range("D5").Value = "'316.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.36%  "

$ws.Range("D6").Value = "'104.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.44%  "

$ws.Range("E7").Value = "  -1.52%  "

$ws.Range("E9").Value = "  -1.09%  "

$ws.Range("D10").Value = "'39.95"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.84%  "

$ws.Range("D11").Value = "'0.0910"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.16%  "

$ws.Range("D12").Value = "'8.48"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.66%  "

$ws.Range("E13").Value = "  +1.29%  "

$ws.Range("D14").Value = "'0.979"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.32%  "

$ws.Range("D15").Value = "'15.42"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.53%  "

$ws.Range("D16").Value = "2.648.77"
$ws.Range("E16").Value = "  -1.94%  "

$ws.Range("D17").Value = "2.293.94"
$ws.Range("E17").Value = "  -1.96%  "

$ws.Range("D18").Value = "41.987.61"

$ws.Range("D19").Value = "'7.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.58%  "

$ws.Range("E20").Value = "  -0.39%  "

$ws.Range("D21").Value = "'73.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.84%  "

$ws.Range("E22").Value = "  -0.73%  "

$ws.Range("D23").Value = "'259.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.25%  "

$ws.Range("E24").Value = "  -0.58%  "

$ws.Range("D25").Value = "'9.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.53%  "

$ws.Range("E26").Value = "  +0.70%  "

$ws.Range("D27").Value = "'10.94"
$ws.Range("D27").Style = "Normal"

$ws.Range("D28").Value = "'23.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.50%  "

$ws.Range("E29").Value = "  +0.70%  "

$ws.Range("D30").Value = "'35.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.42%  "

$ws.Range("D31").Value = "'164.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.03%  "

$ws.Range("D32").Value = "'0.0888"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.01%  "

$ws.Range("D33").Value = "'2.92"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.79%  "

$ws.Range("D34").Value = "'5.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.19%  "

$ws.Range("E35").Value = "  +3.41%  "

$ws.Range("D36").Value = "'0.131"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.37%  "

$ws.Range("D37").Value = "'4.65"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.97%  "

$ws.Range("D38").Value = "'2.94"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +9.96%  "

$ws.Range("E39").Value = "  -2.17%  "

$ws.Range("D40").Value = "'3.67"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.97%  "

$ws.Range("D41").Value = "'101.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +14.72%  "

$ws.Range("E42").Value = "  +0.96%  "

$ws.Range("D43").Value = "'71.24"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.86%  "

$ws.Range("E44").Value = "  -1.85%  "

$ws.Range("D46").Value = "'12.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.16%  "

$ws.Range("D47").Value = "'114.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.84%  "

$ws.Range("D48").Value = "'79.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.66%  "

$ws.Range("E49").Value = "  -0.83%  "

$ws.Range("E50").Value = "  -3.23%  "

$ws.Range("E51").Value = "  +2.21%  "
